$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''37.174.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.49%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.028.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -2.76%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.21%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''226.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -3.27%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.610'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -3.98%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.02%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''55.37'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -4.48%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.383'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -2.42%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.0797'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +2.34%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.105'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -3.55%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''2.327.62'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -2.81%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''14.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -5.62%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''20.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.17%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.746'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -3.72%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''5.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -3.17%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.028.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -2.85%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''37.065.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -1.86%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''6.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +1.34%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''68.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -3.08%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.0₃0837'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +0.41%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''223.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -2.48%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +0.10%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''2.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.45%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''2.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -4.91%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''9.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -2.54%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''167.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -2.02%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -4.59%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''18.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -3.70%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''1.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.51%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.118'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -4.02%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -3.73%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.0611'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -4.00%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''4.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -2.88%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -4.66%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +0.31%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -0.07%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -4.02%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''5.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +0.49%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''1.504.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +3.65%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.0219'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -6.91%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''16.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +0.72%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -1.94%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''95.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -5.76%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.0930'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -3.39%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''1.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -4.64%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -0.23%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -4.28%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -1.01%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -10.16%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''2.214.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -2.83%  '
$ws.Range("E51").Style = "Normal"
